# Auto-generated Excel COM-interop script
# Applies the scheduled-runner price/profit refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 225.23077
$ws.Range("I2").Value = 251.27272
$ws.Range("J2").Value = 82
$ws.Range("K2").Value = 251.27272
$ws.Range("L2").Value = 82
$ws.Range("M2").Value = -138.27272
$ws.Range("N2").Value = -308
$ws.Range("H12").Value = 398.8
$ws.Range("I12").Value = 398.5
$ws.Range("K12").Value = 398.5
$ws.Range("M12").Value = -228.5
$ws.Range("H18").Value = 1491.5834
$ws.Range("J18").Value = 1966.3334
$ws.Range("L18").Value = 1966.3334
$ws.Range("N18").Value = -2534.3334
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H28").Value = 3100.0588
$ws.Range("I28").Value = 2814.4
$ws.Range("J28").Value = 5242.5
$ws.Range("K28").Value = 2814.4
$ws.Range("L28").Value = 5242.5
$ws.Range("M28").Value = -2329.4
$ws.Range("N28").Value = -6212.5
$ws.Range("H33").Value = 473.82608
$ws.Range("I33").Value = 445.35294
$ws.Range("J33").Value = 554.5
$ws.Range("K33").Value = 445.35294
$ws.Range("L33").Value = 554.5
$ws.Range("M33").Value = -216.35294
$ws.Range("N33").Value = -1012.5
$ws.Range("H62").Value = 10105326
$ws.Range("I62").Value = 12349843
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 12349843
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -12349219
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 10105326
$ws.Range("I65").Value = 12349843
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 61749215
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -61746095
$ws.Range("N65").Value = -31240
$ws.Range("H98").Value = 2955.724
$ws.Range("I98").Value = 3317.28
$ws.Range("J98").Value = 696
$ws.Range("K98").Value = 3317.28
$ws.Range("L98").Value = 696
$ws.Range("M98").Value = -1819.28
$ws.Range("N98").Value = -3692
$ws.Range("H100").Value = 1366.5555
$ws.Range("I100").Value = 999.8570999999999
$ws.Range("K100").Value = 999.8570999999999
$ws.Range("M100").Value = -458.8570999999999
$ws.Range("H113").Value = 3312.2856
$ws.Range("I113").Value = 3370
$ws.Range("J113").Value = 3235.3333
$ws.Range("K113").Value = 3370
$ws.Range("L113").Value = 3235.3333
$ws.Range("M113").Value = -116
$ws.Range("N113").Value = -9743.3333
$ws.Range("H116").Value = 3434.7856
$ws.Range("I116").Value = 2900.125
$ws.Range("J116").Value = 4147.6665
$ws.Range("K116").Value = 2900.125
$ws.Range("L116").Value = 4147.6665
$ws.Range("M116").Value = 541.875
$ws.Range("N116").Value = -11031.6665
$ws.Range("H122").Value = 2955.724
$ws.Range("I122").Value = 3317.28
$ws.Range("J122").Value = 696
$ws.Range("K122").Value = 9951.84
$ws.Range("L122").Value = 2088
$ws.Range("M122").Value = -7501.84
$ws.Range("N122").Value = -6988
$ws.Range("H127").Value = 2206.5833
$ws.Range("I127").Value = 1244.75
$ws.Range("K127").Value = 3734.25
$ws.Range("M127").Value = 1225.75
$ws.Range("H138").Value = 436643.2
$ws.Range("J138").Value = 513451.28
$ws.Range("L138").Value = 1540353.84
$ws.Range("N138").Value = -1550633.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5468.6294
$ws.Range("I32").Value = 5332.827
$ws.Range("K32").Value = 5332.827
$ws.Range("M32").Value = -5045.827
$ws.Range("H122").Value = 2573.0557
$ws.Range("I122").Value = 2260.6428
$ws.Range("J122").Value = 3666.5
$ws.Range("K122").Value = 6781.928400000001
$ws.Range("L122").Value = 10999.5
$ws.Range("M122").Value = -4331.928400000001
$ws.Range("N122").Value = -15899.5
$ws.Range("H132").Value = 2624.4614
$ws.Range("I132").Value = 2117.0527
$ws.Range("K132").Value = 6351.158100000001
$ws.Range("M132").Value = -3821.158100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 123.333336
$ws.Range("I11").Value = 123.333336
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 123.333336
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 16.666664
$ws.Range("N11").ClearContents()
$ws.Range("H22").Value = 409.2
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 482
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 482
$ws.Range("M22").Value = -127
$ws.Range("N22").Value = -828
$ws.Range("H94").Value = 14706194
$ws.Range("I94").Value = 22727532
$ws.Range("J94").Value = 406.66666
$ws.Range("K94").Value = 22727532
$ws.Range("L94").Value = 406.66666
$ws.Range("M94").Value = -22727081
$ws.Range("N94").Value = -1308.66666
$ws.Range("H99").Value = 58824610
$ws.Range("I99").Value = 66667744
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 66667744
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = -66666246
$ws.Range("N99").Value = -4096
$ws.Range("H107").Value = 1147.5714
$ws.Range("I107").Value = 736
$ws.Range("J107").Value = 1970.7142
$ws.Range("K107").Value = 736
$ws.Range("L107").Value = 1970.7142
$ws.Range("M107").Value = 1184
$ws.Range("N107").Value = -5810.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 233634.33
$ws.Range("I22").Value = 201
$ws.Range("K22").Value = 201
$ws.Range("M22").Value = 149
$ws.Range("H105").Value = 720.625
$ws.Range("I105").Value = 701.2
$ws.Range("J105").Value = 817.75
$ws.Range("K105").Value = 701.2
$ws.Range("L105").Value = 817.75
$ws.Range("M105").Value = 1045.8
$ws.Range("N105").Value = -4311.75
$ws.Range("H122").Value = 11468.7
$ws.Range("I122").Value = 12520.777
$ws.Range("K122").Value = 37562.331
$ws.Range("M122").Value = -35112.331
$ws.Range("H134").Value = 15626768
$ws.Range("I134").Value = 1759.6428
$ws.Range("J134").Value = 125001820
$ws.Range("K134").Value = 5278.928400000001
$ws.Range("L134").Value = 375005460
$ws.Range("M134").Value = -2743.928400000001
$ws.Range("N134").Value = -375010530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 7755
$ws.Range("I82").Value = 1665.2
$ws.Range("K82").Value = 4995.6
$ws.Range("M82").Value = -4589.6
$ws.Range("H85").Value = 7755
$ws.Range("I85").Value = 1665.2
$ws.Range("K85").Value = 4995.6
$ws.Range("M85").Value = -3591.6
$ws.Range("H122").Value = 1689.0588
$ws.Range("J122").Value = 1741.5
$ws.Range("L122").Value = 15673.5
$ws.Range("N122").Value = -20573.5
$ws.Range("H131").Value = 12987898
$ws.Range("J131").Value = 921.1389
$ws.Range("L131").Value = 2763.4167
$ws.Range("N131").Value = -12843.4167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8764706
$ws.Range("I11").Value = 8076923
$ws.Range("K11").Value = 8076923
$ws.Range("M11").Value = -8076784
$ws.Range("H80").Value = 5931.2
$ws.Range("I80").Value = 6425
$ws.Range("J80").Value = 5602
$ws.Range("K80").Value = 6425
$ws.Range("L80").Value = 5602
$ws.Range("M80").Value = -5427
$ws.Range("N80").Value = -7598
$ws.Range("H83").Value = 5931.2
$ws.Range("I83").Value = 6425
$ws.Range("J83").Value = 5602
$ws.Range("K83").Value = 32125
$ws.Range("L83").Value = 28010
$ws.Range("M83").Value = -27133
$ws.Range("N83").Value = -37994
$ws.Range("H113").Value = 1249.9048
$ws.Range("I113").Value = 1018.2143
$ws.Range("K113").Value = 1018.2143
$ws.Range("M113").Value = 1151.7857
$ws.Range("H122").Value = 2405.889
$ws.Range("I122").Value = 2511.45
$ws.Range("K122").Value = 7534.349999999999
$ws.Range("M122").Value = -5084.349999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H16").Value = 640.2941
$ws.Range("I16").Value = 640.2941
$ws.Range("K16").Value = 640.2941
$ws.Range("M16").Value = -470.2941
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H68").Value = 1878.2858
$ws.Range("I68").Value = 1868.9231
$ws.Range("K68").Value = 1868.9231
$ws.Range("M68").Value = -1119.9231
$ws.Range("H71").Value = 1878.2858
$ws.Range("I71").Value = 1868.9231
$ws.Range("K71").Value = 9344.6155
$ws.Range("M71").Value = -5600.6155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H122").Value = 22728954
$ws.Range("I122").Value = 27779610
$ws.Range("J122").Value = 1005
$ws.Range("K122").Value = 83338830
$ws.Range("L122").Value = 3015
$ws.Range("M122").Value = -83336380
$ws.Range("N122").Value = -7915
$ws.Range("H132").Value = 1457.9062
$ws.Range("I132").Value = 1126.0344
$ws.Range("K132").Value = 3378.1032
$ws.Range("M132").Value = -848.1032
